# Add fake data for database
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix up a couple of existing rows that had data-entry mistakes.
$ws.Cells.Item(2, 4).Value = "Available"
$ws.Cells.Item(3, 2).Value = "SB3"
$ws.Cells.Item(4, 2).Value = "Out"
$ws.Cells.Item(4, 3).Value = "No"
$ws.Cells.Item(4, 4).Value = "Borrowed"

# New fake cup records for the database.
$newRows = @(
    @(989968, "BC10", "No",  "Available"),
    @(989969, "Out",  "No",  "Borrowed"),
    @(989970, "Out",  "Yes", "Sold"),
    @(989971, "BC10", "No",  "Available"),
    @(989972, "BC10", "No",  "Available"),
    @(989973, "BC10", "No",  "Available"),
    @(989974, "BC10", "No",  "Available"),
    @(989975, "BC10", "No",  "Available"),
    @(989976, "SB3",  "No",  "Available"),
    @(989977, "SB3",  "No",  "Available"),
    @(989978, "BC10", "No",  "Available"),
    @(989979, "Out",  "Sold", "Sold")
)

$row = 5
foreach ($rec in $newRows) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $row++
}

$ws.Range("B2:B16").Select() | Out-Null
